$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: "Innleben mit VR..." task ---
$ws.Range("D16").Value = 25
$ws.Range("F16").Value = "x"
$ws.Range("H16").Value = "Alpha Version (noch nicht auslieferbar)"
$ws.Range("H16").WrapText = $true

# --- Row 17: "Verschiedene Moeglichkeiten..." task ---
$ws.Range("D17").Value = 16
$ws.Range("F17").Value = "x"
$ws.Range("H17").Value = "Source Code identifiziert, Darstellung derzeit nicht im Play Mode (Laufzeit) aktualisiert"
$ws.Range("H17").WrapText = $true

# --- Row 18: "UI Buttons in VR hinzugefuegt..." task (text updated) ---
$ws.Range("B18").Value = "UI Buttons in VR hinzugefügt (Akzeptanzkriterien: UI Buttons reagieren auf Handtracking)"
$ws.Range("D18").Value = 7
$ws.Range("F18").Value = "x"
$ws.Range("H18").Value = "Alpha Version (noch nicht auslieferbar)"
$ws.Range("H18").WrapText = $true

# --- Row 19: "Untersuchen von alternativen Workflow..." task (text updated) ---
$ws.Range("B19").Value = "Untersuchen von alternativen Workflow mit GitHub und OneDrive (Akzeptanzkriterium: es wurde ein alternativer Workflow für remote Repositories Untersucht, der größere Dateien als GitHub erlaubt, Entscheidung über zukünftige Verwendung getroffen)"
$ws.Range("D19").Value = 4
$ws.Range("F19").Value = "x"

# --- Totals row: sum of new "Real Effort" column ---
$ws.Range("D27").Formula = "=SUM(D16:D19)"

# --- View state: scroll / zoom / selection as left by the author ---
$ws.Application.ActiveWindow.Zoom = 145
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B19").Select()
